# Add a new "Event " column to the Card18 sheet (right after the "Date" column),
# matching header style of the existing header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card18")

# L1 ("Date") is the last header cell; the new "Event " header goes in M1.
$headerCell = $ws.Range("M1")

# Copy the formatting (bold/centered/bordered) of the existing header cell
# so the new header matches the rest of row 1.
$ws.Range("L1").Copy()
$headerCell.PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Set the new header text (trailing space intentional, matches source data).
$headerCell.Value = "Event "
